# Refresh the cryptocurrency price / 1h-volume-change data pulled from
# coinranking.com for this run of the GitHub Actions scrape.
#
# D (Price) and E (Volume(1h)) columns hold text-formatted numbers (not
# real numeric cells) so values like "139.00" or "  +0.25%  " keep their
# exact original formatting (trailing zeros / padding). Column D values
# that parse as plain numbers need to be forced back to text, otherwise
# the COM layer would silently coerce them to numeric cells and drop
# formatting (e.g. "139.00" -> 139, "218.40" -> 218.4). NumberFormat is
# restored via ClearFormats() right after so no visible/style change is
# left behind - only the cell content changes, exactly like the source
# scrape script would overwrite these cells with fresh text each run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

$ws.Range("D2").Value = "26.223.18"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.661.34"
$ws.Range("E3").Value = "  -0.71%  "
Set-TextValue "D4" "1.004"
$ws.Range("E4").Value = "  +0.16%  "
Set-TextValue "D5" "218.40"
$ws.Range("E5").Value = "  +0.88%  "
Set-TextValue "D6" "0.5225"
$ws.Range("E6").Value = "  -1.51%  "
$ws.Range("E7").Value = "  +0.25%  "
Set-TextValue "D8" "0.2670"
$ws.Range("E8").Value = "  -0.13%  "
Set-TextValue "D9" "0.06324"
$ws.Range("E9").Value = "  -1.08%  "
Set-TextValue "D10" "21.08"
$ws.Range("E10").Value = "  -2.07%  "
Set-TextValue "D11" "0.07716"
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").Value = "1.663.68"
$ws.Range("E12").Value = "  -0.78%  "
Set-TextValue "D13" "4.431"
$ws.Range("E13").Value = "  -1.86%  "
$ws.Range("D14").Value = "1.889.82"
$ws.Range("E14").Value = "  -0.64%  "
Set-TextValue "D15" "0.5468"
$ws.Range("E15").Value = "  -1.73%  "
$ws.Range("D16").Value = "0.0₅8233"
$ws.Range("E16").Value = "  -1.52%  "
Set-TextValue "D17" "64.88"
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "26.256.30"
$ws.Range("E18").Value = "  -0.76%  "
$ws.Range("E19").Value = "  +0.11%  "
Set-TextValue "D20" "4.667"
$ws.Range("E20").Value = "  -2.30%  "
Set-TextValue "D21" "193.42"
$ws.Range("E21").Value = "  -1.03%  "
Set-TextValue "D22" "10.14"
$ws.Range("E22").Value = "  -2.01%  "
Set-TextValue "D23" "6.079"
$ws.Range("E23").Value = "  -3.99%  "
Set-TextValue "D24" "1.007"
$ws.Range("E24").Value = "  +0.40%  "
Set-TextValue "D25" "139.00"
Set-TextValue "D26" "0.1239"
$ws.Range("E26").Value = "  -2.93%  "
Set-TextValue "D27" "7.234"
$ws.Range("E27").Value = "  -2.48%  "
Set-TextValue "D28" "16.18"
$ws.Range("E28").Value = "  -0.90%  "
Set-TextValue "D29" "1.409"
$ws.Range("E29").Value = "  -1.76%  "
$ws.Range("E30").Value = "  -1.96%  "
Set-TextValue "D31" "1.282"
$ws.Range("E31").Value = "  +0.59%  "
Set-TextValue "D32" "3.682"
$ws.Range("E32").Value = "  +1.97%  "
Set-TextValue "D33" "3.316"
$ws.Range("E33").Value = "  -3.74%  "
Set-TextValue "D34" "1.634"
$ws.Range("E34").Value = "  -3.40%  "
Set-TextValue "D35" "0.9801"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-TextValue "D36" "2.415"
$ws.Range("E36").Value = "  -0.43%  "
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D37" "2.782"
$ws.Range("E37").Value = "  +0.16%  "
Set-TextValue "D38" "0.5880"
$ws.Range("E38").Value = "  +3.02%  "
Set-TextValue "D39" "0.01593"
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("E40").Value = "  -1.11%  "
Set-TextValue "D41" "0.8596"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "1.032.51"
$ws.Range("E43").Value = "  -3.37%  "
Set-TextValue "D44" "99.59"
$ws.Range("E44").Value = "  -0.41%  "
$ws.Range("D45").Value = "1.804.02"
$ws.Range("E45").Value = "  -1.02%  "
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").Value = "  +2.21%  "
Set-TextValue "D47" "57.17"
$ws.Range("E47").Value = "  +0.26%  "
$ws.Range("E48").Value = "  +0.02%  "
Set-TextValue "D49" "8.095"
$ws.Range("E49").Value = "  -0.30%  "
Set-TextValue "D50" "0.05186"
$ws.Range("E50").Value = "  -0.43%  "
Set-TextValue "D51" "1.470"
$ws.Range("E51").Value = "  +0.34%  "
